$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values (column D) are stored as plain text in the sheet (e.g. two-dot
# numbers like "29.857.74", or numbers that would otherwise be auto-converted
# to a float/scientific by Excel, e.g. "0.992" -> 0.99199999999999999). Force
# text interpretation via NumberFormat="@" then restore the cell to the default
# "Normal" style so no stray number-format style lingers on the cell.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "29.857.74"
$ws.Range("E2").Value = "  +1.17%  "

# Row 3
Set-TextValue "D3" "1.618.08"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4
Set-TextValue "D4" "0.992"
$ws.Range("E4").Value = "  -0.85%  "

# Row 5
Set-TextValue "D5" "213.52"
$ws.Range("E5").Value = "  +0.31%  "

# Row 6
$ws.Range("E6").Value = "  -0.69%  "

# Row 7
Set-TextValue "D7" "0.991"
$ws.Range("E7").Value = "  -0.91%  "

# Row 8
Set-TextValue "D8" "29.23"
$ws.Range("E8").Value = "  +8.99%  "

# Row 9
$ws.Range("E9").Value = "  +3.25%  "

# Row 10
$ws.Range("E10").Value = "  +0.84%  "

# Row 11
Set-TextValue "D11" "0.0911"
$ws.Range("E11").Value = "  -0.04%  "

# Row 12
Set-TextValue "D12" "1.851.50"
$ws.Range("E12").Value = "  +0.70%  "

# Row 13
Set-TextValue "D13" "1.615.09"
$ws.Range("E13").Value = "  +0.25%  "

# Row 14
$ws.Range("E14").Value = "  +5.67%  "

# Row 15
$ws.Range("E15").Value = "  +4.95%  "

# Row 16
Set-TextValue "D16" "29.886.11"
$ws.Range("E16").Value = "  +1.21%  "

# Row 17
Set-TextValue "D17" "8.84"
$ws.Range("E17").Value = "  +15.64%  "

# Row 18
Set-TextValue "D18" "64.37"
$ws.Range("E18").Value = "  +1.51%  "

# Row 19
Set-TextValue "D19" "241.10"
$ws.Range("E19").Value = "  -0.06%  "

# Row 20
Set-TextValue "D20" "0.0₃0708"
$ws.Range("E20").Value = "  +2.52%  "

# Row 21
Set-TextValue "D21" "0.993"
$ws.Range("E21").Value = "  -0.72%  "

# Row 22
Set-TextValue "D22" "4.10"
$ws.Range("E22").Value = "  +2.34%  "

# Row 23
Set-TextValue "D23" "9.61"
$ws.Range("E23").Value = "  +4.33%  "

# Row 24
$ws.Range("E24").Value = "  +0.50%  "

# Row 25
Set-TextValue "D25" "155.11"
$ws.Range("E25").Value = "  +0.28%  "

# Row 26
Set-TextValue "D26" "15.59"
$ws.Range("E26").Value = "  +2.13%  "

# Row 27
Set-TextValue "D27" "0.110"
$ws.Range("E27").Value = "  +1.03%  "

# Row 28
Set-TextValue "D28" "6.58"
$ws.Range("E28").Value = "  +3.08%  "

# Row 29
Set-TextValue "D29" "0.993"
$ws.Range("E29").Value = "  -0.72%  "

# Row 30
$ws.Range("E30").Value = "  +3.09%  "

# Row 31
Set-TextValue "D31" "1.13"
$ws.Range("E31").Value = "  +5.91%  "

# Row 32
$ws.Range("E32").Value = "  +3.30%  "

# Row 33
$ws.Range("E33").Value = "  +3.38%  "

# Row 34
Set-TextValue "D34" "1.415.97"
$ws.Range("E34").Value = "  +0.12%  "

# Row 35
$ws.Range("E35").Value = "  +6.41%  "

# Row 36
$ws.Range("E36").Value = "  +0.15%  "

# Row 37
Set-TextValue "D37" "2.85"
$ws.Range("E37").Value = "  +1.19%  "

# Row 38
Set-TextValue "D38" "2.29"
$ws.Range("E38").Value = "  -0.77%  "

# Row 39
$ws.Range("E39").Value = "  +2.24%  "

# Row 40
Set-TextValue "D40" "0.557"
$ws.Range("E40").Value = "  +3.70%  "

# Row 41
$ws.Range("E41").Value = "  +3.07%  "

# Row 42
$ws.Range("E42").Value = "  +0.27%  "

# Row 43
Set-TextValue "D43" "0.826"
$ws.Range("E43").Value = "  +3.49%  "

# Row 44
Set-TextValue "D44" "53.68"
$ws.Range("E44").Value = "  +2.46%  "

# Row 45
Set-TextValue "D45" "69.25"
$ws.Range("E45").Value = "  +5.25%  "

# Row 46
$ws.Range("E46").Value = "  +18.64%  "

# Row 47
$ws.Range("E47").Value = "  -0.96%  "

# Row 48
$ws.Range("E48").Value = "  +2.82%  "

# Row 49
Set-TextValue "D49" "1.759.96"
$ws.Range("E49").Value = "  +0.66%  "

# Row 50
Set-TextValue "D50" "88.25"
$ws.Range("E50").Value = "  +1.76%  "

# Row 51
Set-TextValue "D51" "0.0532"
$ws.Range("E51").Value = "  +1.87%  "
